$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D16").Value = "2016-03-08 02:39:20"
$wsZh.Range("G16").Value = "2016-03-08 02:40:00"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D16").Value = "2016-03-08 02:39:28"
$wsDe.Range("G16").Value = "2016-03-08 02:40:15"
